$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.810.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.572.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.576"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.24%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.549"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0813"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.25%  "

# Row 13
$ws.Range("E13").Value = "  +6.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.568.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.886"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "42.838.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.97%  "

# Row 18
$ws.Range("B18").Value = "InternetComputer(DFINITY)"
$ws.Range("C18").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.83%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("E20").Value = "  -1.30%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.99%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "254.22"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("E24").Value = "  -3.97%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.27%  "

# Row 29
$ws.Range("E29").Value = "  -2.82%  "

# Row 30
$ws.Range("E30").Value = "  -2.20%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.87%  "

# Row 32
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.78%  "

# Row 33
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.41%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0803"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.35%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.06%  "

# Row 37
$ws.Range("E37").Value = "  -2.50%  "

# Row 38
$ws.Range("E38").Value = "  -0.95%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.85%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +29.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.91%  "

# Row 42
$ws.Range("E42").Value = "  +0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.084.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.15%  "

# Row 45
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("E46").Value = "  +0.82%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.30%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.49%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.822.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.192"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.13%  "
